# feat: add 2022-Q4 data
#
# Before:  Sheet1 "总计" (summary) + Sheet2 "2022-Q3" (fund snapshot)
# After:   Sheet1 "总计" (summary, +1 row) + Sheet2 "2022-Q4" (new fund
#          snapshot, same underlying sheet/id as the old "2022-Q3") +
#          Sheet3 "2022-Q3" (new sheet holding what used to be in sheet2)

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)

# 1. Duplicate the existing "2022-Q3" sheet so its current data survives
#    under its own tab, placed right after it.
$q3.Copy($null, $q3)
$newQ3 = $wb.Worksheets.Item(3)

# 2. Turn the original sheet into the new "2022-Q4" snapshot (rename it
#    out of the way first so the duplicate can reclaim the "2022-Q3" name).
$q3.Name = "2022-Q4"
$newQ3.Name = "2022-Q3"

# Overwrite the fund snapshot figures with the 2022-Q4 numbers. Use a
# leading quote so the numeric-looking text ("0.24", "68.42", ...) is
# stored as text, matching the source column's type - then strip the
# formatting PasteSpecial left behind by pulling a plain (unstyled)
# neighbour cell's format over it.
$q3.Range("D2").Value = "'0.24"
$q3.Range("E2").Value = "'68.42"
$q3.Range("F2").Value = "'3.94"
$q3.Range("G2").Value = "'0.0095"
$q3.Range("H2").Value = 9

$q3.Range("C2").Copy()
$q3.Range("D2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The new snapshot's header row + index cell pick up the "总计" sheet's
# header style rather than the old tab's.
$summary.Range("A2").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Record the new quarter in the "总计" summary sheet: rename the
#    existing data row to 2022-Q4 and append a 2022-Q3 row below it.
$summary.Range("B2").Value = "2022-Q4"

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.01

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
